$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.099.65'
$ws.Range("E2").Value = '  +0.67%  '
$ws.Range("D3").Value = '2.569.60'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '505.16'
$ws.Range("E5").Value = '  -0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -5.18%  '
$ws.Range("D9").Value = '2.572.36'
$ws.Range("E9").Value = '  -0.21%  '
$ws.Range("D10").Value = '6.56'
$ws.Range("E10").Value = '  +7.26%  '
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("E12").Value = '  +1.59%  '
$ws.Range("E13").Value = '  +1.22%  '
$ws.Range("D14").Value = '3.021.36'
$ws.Range("E14").Value = '  +0.05%  '
$ws.Range("D15").Value = '60.170.95'
$ws.Range("E15").Value = '  +1.11%  '
$ws.Range("D16").Value = '21.49'
$ws.Range("E16").Value = '  -1.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000139'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.55%  '
$ws.Range("D18").Value = '2.570.34'
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '344.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.68%  '
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("E23").Value = '  -0.43%  '
$ws.Range("E24").Value = '  -0.69%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("E26").Value = '  +0.36%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("E28").Value = '  +0.90%  '
$ws.Range("D29").Value = '7.38'
$ws.Range("E29").Value = '  +0.80%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").Value = '19.28'
$ws.Range("E31").Value = '  -0.77%  '
$ws.Range("D32").Value = '153.08'
$ws.Range("E33").Value = '  -1.07%  '
$ws.Range("E34").Value = '  +3.08%  '
$ws.Range("E35").Value = '  +1.54%  '
$ws.Range("E36").Value = '  -1.08%  '
$ws.Range("D37").Value = '0.848'
$ws.Range("E37").Value = '  +8.29%  '
$ws.Range("D38").Value = '0.846'
$ws.Range("E38").Value = '  -1.39%  '
$ws.Range("E39").Value = '  +1.96%  '
$ws.Range("D40").Value = '36.07'
$ws.Range("E40").Value = '  +2.40%  '
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").Value = '293.05'
$ws.Range("E42").Value = '  -3.04%  '
$ws.Range("E43").Value = '  -1.90%  '
$ws.Range("E44").Value = '  -2.54%  '
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("E46").Value = '  -3.05%  '
$ws.Range("D47").Value = '19.68'
$ws.Range("E47").Value = '  +1.87%  '
$ws.Range("E48").Value = '  -2.40%  '
$ws.Range("E49").Value = '  -1.71%  '
$ws.Range("D50").Value = '10.31'
$ws.Range("E50").Value = '  +0.45%  '
$ws.Range("D51").Value = '1.991.54'
$ws.Range("E51").Value = '  +0.20%  '
